$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.840.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +4.70%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'3.059.74"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +2.36%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'576.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +2.58%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'142.90"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +3.68%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'3.056.03"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +2.40%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  +0.77%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.139"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +4.46%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'5.47"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +11.87%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.464"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +1.59%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +4.00%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'34.86"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +3.49%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  -0.14%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.565.54"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +2.44%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'7.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +3.07%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'3.057.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +2.46%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'61.746.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +4.69%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'449.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +5.91%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'13.95"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +3.03%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.732"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +2.96%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'7.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +2.78%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'13.69"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +1.90%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'81.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +1.79%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +0.06%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'2.24"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +4.97%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +0.01%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  +3.96%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'8.04"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +3.14%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'6.63"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +8.45%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +3.72%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'0.107"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +7.91%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'0.0₃0809"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +4.05%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  +2.08%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'6.07"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +5.69%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'2.20"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +5.83%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'50.14"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +2.31%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'2.95"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +5.78%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'8.83"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +2.19%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'416.02"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +3.35%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.0367"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +5.22%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'2.772.38"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.98%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  +0.77%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.265"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +7.68%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'36.61"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +12.38%  "
$ws.Range("E46").ClearFormats()
$ws.Range("B47").Value = "'Fetch.AI"
$ws.Range("B47").ClearFormats()
$ws.Range("C47").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C47").ClearFormats()
$ws.Range("D47").Value = "'2.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +3.48%  "
$ws.Range("E47").ClearFormats()
$ws.Range("B48").Value = "'USDe"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +0.02%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'123.07"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -1.61%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  +1.54%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'24.08"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +2.98%  "
$ws.Range("E51").ClearFormats()
